$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before HB (column 210). This shifts the existing
# "nom" (HB) and "url_produit" (HC) columns one position to the right
# (to HC and HD respectively), matching the dimension change from
# A1:HC210 to A1:HD210.
$ws.Columns("HB:HB").Insert()

# New header cell for the newly inserted timestamp column, with the same
# (bold/centered/bordered) header style used by the other snapshot columns.
$ws.Range("HB1").Value = "2026-02-06 17:33:41"
$ws.Range("HA1").Copy()
$ws.Range("HB1").PasteSpecial(-4122)

# For the rows that already had a price captured in the previous last
# snapshot column (HA, now still column 209), duplicate that value into
# the freshly inserted column (HB, column 210) - this represents the new
# price snapshot being identical to the last recorded one.
for ($r = 2; $r -le 80; $r++) {
    $srcCell = $ws.Cells.Item($r, 209)
    $dstCell = $ws.Cells.Item($r, 210)
    $dstCell.Value = $srcCell.Value2
}

# Rows 81-210 had no value in HA, so the newly inserted HB cell stays
# empty for those rows (default state after the column insert).
